$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old growth-group data block (rows 30-71: previously the 0.4/0.8/1.6/3.2
# groups followed by two stray Growth_0.3 rows at the very end).
$ws.Rows("30:71").Delete()

# Rebuild the data block starting at row 30. The 0.3 growth-rate group is now fully
# populated (NoMR_1..4, MR_1..5) and sits in its correct sorted position ahead of the
# 0.4/0.8/1.6/3.2 groups, which are reproduced unchanged.
$groups = @(
    @(0.3, @("NoMR_1", "NoMR_2", "NoMR_3", "NoMR_4", "MR_1", "MR_2", "MR_3", "MR_4", "MR_5")),
    @(0.4, @("NoMR_1", "NoMR_2", "NoMR_3", "NoMR_4", "NoMR_5", "MR_1", "MR_2", "MR_3", "MR_4", "MR_5")),
    @(0.8, @("NoMR_1", "NoMR_2", "NoMR_3", "NoMR_4", "NoMR_5", "MR_1", "MR_2", "MR_3", "MR_4", "MR_5")),
    @(1.6, @("NoMR_1", "NoMR_2", "NoMR_3", "NoMR_4", "NoMR_5", "MR_1", "MR_2", "MR_3", "MR_4", "MR_5")),
    @(3.2, @("NoMR_1", "NoMR_2", "NoMR_3", "NoMR_4", "NoMR_5", "MR_1", "MR_2", "MR_3", "MR_4", "MR_5"))
)

$row = 30
foreach ($grp in $groups) {
    $g = $grp[0]
    $labels = $grp[1]
    foreach ($lbl in $labels) {
        $isMR = $lbl.StartsWith("MR_")
        $ws.Cells.Item($row, 1).Value = "Exp1_GrowthAndMR"
        $ws.Cells.Item($row, 2).Value = "Growth_${g}_${lbl}"
        $ws.Cells.Item($row, 3).Value = $g
        $ws.Cells.Item($row, 4).Value = $isMR
        $ws.Cells.Item($row, 5).Value = $g * 0.75
        $ws.Cells.Item($row, 6).Value = $g * 1.25
        $row = $row + 1
    }
}

$lastRow = $row - 1

# Re-apply the sort (by growth rate, then by MR flag) so the worksheet's remembered
# AutoFilter/sort state covers the new A2:F78 extent instead of the old A2:D80.
$sortRange = $ws.Range("A2:F" + $lastRow)
$s = $ws.Sort
$s.SortFields.Clear()
[void]$s.SortFields.Add($ws.Range("C2:C" + $lastRow))
[void]$s.SortFields.Add($ws.Range("D2:D" + $lastRow))
$s.SetRange($sortRange)
$s.Header = 0
$s.Apply()

# Restore the active cell/selection to where the author last left off.
[void]$ws.Range("H24").Select()
